$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 4).Value = 150.0
$ws.Cells.Item(8, 5).Value = -1.0
$ws.Cells.Item(8, 6).Value = 9.0
$ws.Cells.Item(8, 7).Value = 34.0
$ws.Cells.Item(8, 8).Value = 13.0
$ws.Cells.Item(8, 9).Value = 11.0
$ws.Cells.Item(8, 10).Value = 12.0
$ws.Cells.Item(8, 11).Value = 16.0
$ws.Cells.Item(8, 12).Value = 19.0
$ws.Cells.Item(8, 13).Value = 20.0
$ws.Cells.Item(8, 14).Value = 21.0
$ws.Cells.Item(8, 15).Value = 26.0
$ws.Cells.Item(8, 16).Value = 29.0
$ws.Cells.Item(8, 17).Value = 28.0
$ws.Cells.Item(8, 18).Value = 27.0
$ws.Cells.Item(8, 19).Value = 18.0
$ws.Cells.Item(8, 20).Value = 15.0
$ws.Cells.Item(8, 21).Value = 14.0
$ws.Cells.Item(8, 22).Value = 17.0
$ws.Cells.Item(8, 23).Value = 33.0
$ws.Cells.Item(8, 24).Value = 31.0
$ws.Cells.Item(8, 25).Value = 30.0
$ws.Cells.Item(8, 26).Value = 36.0
$ws.Cells.Item(8, 27).Value = 42.0
$ws.Cells.Item(8, 28).Value = 41.0
$ws.Cells.Item(8, 29).Value = 45.0
$ws.Cells.Item(8, 30).Value = 50.0
$ws.Cells.Item(8, 31).Value = 49.0
$ws.Cells.Item(8, 32).Value = 57.0
$ws.Cells.Item(8, 33).Value = 59.0
$ws.Cells.Item(8, 34).Value = 48.0
$ws.Cells.Item(8, 35).Value = 47.0
$ws.Cells.Item(8, 36).Value = 44.0
$ws.Cells.Item(8, 37).Value = 10.0
$ws.Cells.Item(8, 38).Value = 8.0
$ws.Cells.Item(8, 39).Value = 7.0
$ws.Cells.Item(8, 40).Value = 6.0
$ws.Cells.Item(8, 41).Value = 5.0
$ws.Cells.Item(8, 42).Value = 39.0
$ws.Cells.Item(8, 43).Value = 40.0
$ws.Cells.Item(8, 44).Value = 3.0
$ws.Cells.Item(8, 45).Value = 4.0
$ws.Cells.Item(8, 46).Value = -1.0

# Row 9
$ws.Cells.Item(9, 5).Value = -1.0

# Row 10
$ws.Cells.Item(10, 4).Value = 121.0
$ws.Cells.Item(10, 5).Value = -1.0
$ws.Cells.Item(10, 16).Value = 42.0
$ws.Cells.Item(10, 17).Value = 41.0
$ws.Cells.Item(10, 18).Value = 44.0
$ws.Cells.Item(10, 19).Value = 45.0
$ws.Cells.Item(10, 20).Value = 50.0
$ws.Cells.Item(10, 21).Value = 49.0
$ws.Cells.Item(10, 22).Value = 57.0
$ws.Cells.Item(10, 23).Value = 59.0
$ws.Cells.Item(10, 24).Value = 48.0
$ws.Cells.Item(10, 25).Value = 47.0
$ws.Cells.Item(10, 26).Value = 39.0
$ws.Cells.Item(10, 27).Value = 40.0
$ws.Cells.Item(10, 28).Value = -1.0

# Row 11
$ws.Cells.Item(11, 5).Value = -1.0

# Row 12
$ws.Cells.Item(12, 4).Value = 116.0
$ws.Cells.Item(12, 5).Value = -1.0
$ws.Cells.Item(12, 6).Value = 40.0
$ws.Cells.Item(12, 7).Value = 47.0
$ws.Cells.Item(12, 8).Value = 48.0
$ws.Cells.Item(12, 9).Value = 59.0
$ws.Cells.Item(12, 10).Value = 57.0
$ws.Cells.Item(12, 11).Value = 50.0
$ws.Cells.Item(12, 12).Value = 49.0
$ws.Cells.Item(12, 13).Value = 44.0
$ws.Cells.Item(12, 14).Value = 45.0
$ws.Cells.Item(12, 15).Value = 30.0
$ws.Cells.Item(12, 16).Value = 28.0
$ws.Cells.Item(12, 17).Value = 29.0
$ws.Cells.Item(12, 18).Value = 27.0
$ws.Cells.Item(12, 19).Value = 26.0
$ws.Cells.Item(12, 20).Value = 18.0
$ws.Cells.Item(12, 21).Value = 31.0
$ws.Cells.Item(12, 22).Value = 34.0
$ws.Cells.Item(12, 23).Value = 36.0
$ws.Cells.Item(12, 24).Value = 42.0
$ws.Cells.Item(12, 25).Value = 41.0
$ws.Cells.Item(12, 26).Value = -1.0

# Row 13
$ws.Cells.Item(13, 5).Value = -1.0
